# Title-case a subset of "Type" values (column C) in the Topsongs sheet.
# Genres that were lowercase single/compound words describing "pop" sub-genres
# get their first letter capitalised; compound genres already containing an
# ampersand/abbreviation (EDM pop, Orchestral pop, R&B pop) are left as-is.

$map = @{
  "pop"          = "Pop"
  "hiphop"       = "Hiphop"
  "acoustic pop" = "Acoustic pop"
  "disco"        = "Disco"
  "folk"         = "Folk"
  "dance pop"    = "Dance pop"
  "soul pop"     = "Soul pop"
  "pop rock"     = "Pop rock"
  "indie pop"    = "Indie pop"
  "country pop"  = "Country pop"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Topsongs")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 132 }

for ($r = 2; $r -le $lastRow; $r++) {
  $cell = $ws.Cells.Item($r, 3)
  $val = $cell.Text
  if ($map.ContainsKey($val)) {
    $cell.Value = $map[$val]
  }
}

# Reflect the author's final on-screen selection/scroll position from the
# saved file (row ~67 in view, active cell C77).
$ws.Range("C77").Select()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 1
